{"js": "// Office.js (Word JavaScript API) script\n// Change applied (matches the commit \"feat: Add progress #38 'Learn Typescript'\"):\n//   1. Split the last paragraph's run \":, number\" into two runs \":,\" / \" number\"\n//      separated by proofErr gramStart/gramEnd markers (Word's grammar-check\n//      artifacts around the colon).\n//   2. Append a new blank paragraph.\n//   3. Append a new quiz question paragraph (list level 0):\n//        \"Which function call will result in a compilation error?\"\n//   4. Append a new quiz answer paragraph (list level 1):\n//        \"greet('Hi','Bill', 'Gates');\"\n//      with a spellStart/spellEnd pair wrapping \"Hi','Bill\" (Word's spell-check\n//      artifact around that token).\n\nconst runProps =\n  '<w:rPr><w:rFonts w:ascii=\"AppleSystemUIFont\" w:hAnsi=\"AppleSystemUIFont\" ' +\n  'w:cs=\"AppleSystemUIFont\"/><w:lang w:val=\"en-US\"/></w:rPr>';\n\n// Wrap a <w:body> fragment in the minimal OOXML package insertOoxml expects.\nfunction wrapPackage(bodyXml) {\n  return (\n    '<?xml version=\"1.0\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" ' +\n    'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    bodyXml +\n    '<w:sectPr/></w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\n// Locate the last paragraph in the document body (the \":, number\" quiz answer).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst lastIndex = paragraphs.items.length - 1;\nconst lastParagraph = paragraphs.items[lastIndex];\nlastParagraph.load(\"text\");\nawait context.sync();\n\n// --- Step 1: split \":, number\" into \":,\" + \" number\" with proofErr marks ---\nconst splitFragment =\n  '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  \"<w:r>\" + runProps + \"<w:t>:,</w:t></w:r>\" +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  \"<w:r>\" + runProps + '<w:t xml:space=\"preserve\"> number</w:t></w:r>' +\n  \"</w:p>\";\n\n// Insert the new runs right before the paragraph's existing content; when the\n// target range is collapsed at the paragraph start, InsertLocation.before\n// merges the incoming runs into the paragraph instead of creating a new one.\nconst startRange = lastParagraph.getRange(\"Start\");\nstartRange.insertOoxml(wrapPackage(splitFragment), Word.InsertLocation.before);\nawait context.sync();\n\n// Now the paragraph text reads \":, number:, number\" (new runs + old run).\n// Remove the second (original) occurrence, leaving only the new split runs.\nconst duplicateMatches = lastParagraph.search(\":, number\", { matchCase: true });\nduplicateMatches.load(\"text\");\nawait context.sync();\nduplicateMatches.items[1].insertText(\"\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Step 2: append a blank paragraph + the two new quiz paragraphs ---\nconst emptyParagraphXml =\n  '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"></w:p>';\n\nconst questionParagraphXml =\n  '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/>' +\n  '<w:numId w:val=\"11\"/></w:numPr></w:pPr>' +\n  \"<w:r>\" + runProps + \"<w:t>Which function call will result in a compilation error?</w:t></w:r>\" +\n  \"</w:p>\";\n\nconst answerParagraphXml =\n  '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"1\"/>' +\n  '<w:numId w:val=\"11\"/></w:numPr></w:pPr>' +\n  \"<w:r>\" + runProps + \"<w:t>greet('</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  \"<w:r>\" + runProps + \"<w:t>Hi','Bill</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  \"<w:r>\" + runProps + \"<w:t>', 'Gates');</w:t></w:r>\" +\n  \"</w:p>\";\n\nconst endRange = lastParagraph.getRange(\"End\");\nendRange.insertOoxml(\n  wrapPackage(emptyParagraphXml + questionParagraphXml + answerParagraphXml),\n  Word.InsertLocation.after\n);\nawait context.sync();\n", "ps1": "# Word COM interop script\n# Change applied (matches the commit \"feat: Add progress #38 'Learn Typescript'\"):\n#   1. Split the last paragraph's run \":, number\" into two runs \":,\" / \" number\"\n#      separated by proofErr gramStart/gramEnd markers (Word's grammar-check\n#      artifacts around the colon).\n#   2. Append a new blank paragraph.\n#   3. Append a new quiz question paragraph (list level 0):\n#        \"Which function call will result in a compilation error?\"\n#   4. Append a new quiz answer paragraph (list level 1):\n#        \"greet('Hi','Bill', 'Gates');\"\n#      with a spellStart/spellEnd pair wrapping \"Hi','Bill\" (Word's spell-check\n#      artifact around that token).\n\n$d = $word.ActiveDocument\n\n$rpr = \"<w:rPr><w:rFonts w:ascii='AppleSystemUIFont' w:hAnsi='AppleSystemUIFont' w:cs='AppleSystemUIFont'/><w:lang w:val='en-US'/></w:rPr>\"\n$wns = \"xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'\"\n\n# --- Step 1: split the last paragraph's text into two runs with proofErr marks ---\n$paras = $d.Paragraphs\n$count = $paras.Count\n$last = $paras.Item($count)\n\n$r = $last.Range\n$r.MoveEnd(1, -1)      # exclude the paragraph mark from the range\n$r.Text = \"\"           # clear the existing \":, number\" run, keep pPr/list formatting\n\n$pStart = $last.Range.Start\n$insertionPoint = $d.Range($pStart, $pStart)\n$splitXml = \"<w:p $wns>\" + `\n  \"<w:proofErr w:type='gramStart'/>\" + `\n  \"<w:r>$rpr<w:t>:,</w:t></w:r>\" + `\n  \"<w:proofErr w:type='gramEnd'/>\" + `\n  \"<w:r>$rpr<w:t xml:space='preserve'> number</w:t></w:r>\" + `\n  \"</w:p>\"\n$insertionPoint.InsertXML($splitXml)\n\n# --- Step 2: append the blank paragraph + the two new quiz paragraphs ---\n$content = $d.Content\n$endPoint = $d.Range($content.End, $content.End)\n$appendXml = \"<w:p $wns></w:p>\" + `\n  \"<w:p $wns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='11'/></w:numPr></w:pPr>\" + `\n  \"<w:r>$rpr<w:t>Which function call will result in a compilation error?</w:t></w:r></w:p>\" + `\n  \"<w:p $wns><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='11'/></w:numPr></w:pPr>\" + `\n  \"<w:r>$rpr<w:t>greet('</w:t></w:r>\" + `\n  \"<w:proofErr w:type='spellStart'/>\" + `\n  \"<w:r>$rpr<w:t>Hi','Bill</w:t></w:r>\" + `\n  \"<w:proofErr w:type='spellEnd'/>\" + `\n  \"<w:r>$rpr<w:t>', 'Gates');</w:t></w:r></w:p>\"\n$endPoint.InsertXML($appendXml)\n"}
